$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 (shifts existing rows 6-14 down to 7-15).
# Excel inherits formatting from the row above (row 5) automatically.
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = "Growth Stage name"
$ws.Range("B6").Value = "whole plant flowering stage"
$ws.Range("D6").Value = "use PO name from www.plantontology.org"

# Update selection to A7 as indicated in the diff
$ws.Range("A7").Select()
